# Edit cars_missing_rownames worksheet:
# - unmerge A6:A7 so each cylinder-group row gets its own row/value
# - shift the row5/row6/row7 "N / Mean / SD" statistics so row5 now
#   holds the values that used to live in row6, row6 now holds what
#   used to live in row5, and row7's cylinder count cell (A7) gets the
#   value that used to be shown (merged) in A6
# - restyle A6/A7 to match the plain "cylinder count" cell style (A8)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unmerge the A6:A7 block; the cylinder counts will now be shown on two
# separate rows instead of one merged cell.
$ws.Range("A6:A7").UnMerge()

# --- Row 5 values (N / Mean HP / SD HP / Mean WT / SD WT) ---
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 110
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 2.7475
$ws.Range("G5").Value = 0.1803122292025695

# --- Row 6 values ---
$ws.Range("A6").Value = 4
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 83.33333333333333
$ws.Range("E6").Value = 18.50225211517056
$ws.Range("F6").Value = 2.886666666666667
$ws.Range("G6").Value = 0.4911551010967242

# --- Row 7 value ---
$ws.Range("A7").Value = 6

# Match A6/A7 formatting (border/fill/font) to the plain cylinder-count
# cell style already used by A5/A8, since the cells are no longer part
# of a vertically-merged, top-aligned block.
$ws.Range("A8").Copy()
$ws.Range("A6:A7").PasteSpecial(-4122)
$excel.CutCopyMode = 0
